# Update the small Bahnar/English word-pair table on Sheet1.
# Row 2 keeps "hai" -> "two" (N) but the trailing space on "two " is removed.
# Row 3 is changed from "vịt" -> "ara" (N) to "một" -> "one" (N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "hai"
$ws.Range("B2").Value = "two"
$ws.Range("C2").Value = "N"

$ws.Range("A3").Value = "một"
$ws.Range("B3").Value = "one"
$ws.Range("C3").Value = "N"

# Match the final selection recorded in the sheet view.
$ws.Range("C3").Select()
